$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D (removes the "Combined Accuracy" header and all its data)
$ws.Range("D1:D21").Delete()

# Updated Digit Accuracy (B) and Author Accuracy (C) values for the three-model run
$values = @(
    @(39, 23.83333333333333),
    @(58.66666666666666, 12.66666666666667),
    @(77.83333333333333, 20.66666666666667),
    @(86.16666666666667, 22.5),
    @(88.83333333333333, 22),
    @(89, 24.66666666666667),
    @(90.83333333333333, 21.66666666666667),
    @(92.16666666666667, 22.83333333333333),
    @(91.83333333333333, 22.33333333333333),
    @(91.83333333333333, 22.16666666666667),
    @(92, 22.16666666666667),
    @(92.33333333333333, 22.66666666666667),
    @(92.33333333333333, 21.83333333333333),
    @(92.33333333333333, 22.33333333333333),
    @(92.16666666666667, 21.66666666666667),
    @(92.33333333333333, 21.83333333333333),
    @(92.33333333333333, 21.83333333333333),
    @(92.5, 21.83333333333333),
    @(92.5, 21.66666666666667),
    @(92.5, 21.66666666666667)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
}
